$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pepID")

$ws.Range("A29").Value = "PEP_ID-2009652"
$ws.Range("A30").Value = "PEP_ID-2009655"
$ws.Range("A31").Value = "PEP_ID-2009656"
